$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ----------------------------------------------------------------------
# Row 12 (1-indexed): "Uniforms for islamiyya/ shoes" | "₦ 3500/ 1600"
#   -> drop the "/ shoes" suffix (leave a trailing space) and trim the
#      amount down to the uniform-only figure "3500".
# ----------------------------------------------------------------------
$uniformLabel = $t.Cell(12, 3)
$uniformLabel.Range.Find.Execute("/ shoes", $false, $false, $false, $false, $false, `
    $true, 1, $false, " ", 2)

$uniformAmount = $t.Cell(12, 4)
$uniformAmount.Range.Find.Execute("3500/ 1600", $false, $false, $false, $false, $false, `
    $true, 1, $false, "3500", 2)

# ----------------------------------------------------------------------
# Row 13 (1-indexed) used to hold the "Payment of loan" entry. It now
# becomes the "Shoes for islamiyya" entry (the shoe cost that was
# trimmed out of row 12 above), with the amount updated to 1600.
# ----------------------------------------------------------------------
$shoesLabel = $t.Cell(13, 3)
$shoesLabel.Range.Find.Execute(" Payment of loan", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Shoes for islamiyya", 2)

$shoesAmount = $t.Cell(13, 4)
$shoesAmount.Range.Find.Execute("1000", $false, $false, $false, $false, $false, `
    $true, 1, $false, "1600", 2)

# ----------------------------------------------------------------------
# The original "Payment of loan" / "₦ 1000" entry still belongs in the
# table -- re-insert it as its own row directly below, ahead of the
# "Wedding contribution" row.
# ----------------------------------------------------------------------
$weddingRow = $t.Rows.Item(14)
$loanRow = $t.Rows.Add($weddingRow)
$loanRowIndex = $loanRow.Index
$t.Cell($loanRowIndex, 3).Range.Text = " Payment of loan"
$t.Cell($loanRowIndex, 4).Range.Text = "₦ 1000"

# ----------------------------------------------------------------------
# Wedding contribution row (now pushed down by one): tidy the two-run
# label/amount into single runs with identical text.
# ----------------------------------------------------------------------
$weddingRowIndex = $loanRowIndex + 1
$weddingLabel = $t.Cell($weddingRowIndex, 3)
$weddingLabel.Range.Find.Execute(" Wedding contribution", $false, $false, $false, $false, `
    $false, $true, 1, $false, " Wedding contribution", 2)

$weddingAmount = $t.Cell($weddingRowIndex, 4)
$weddingAmount.Range.Find.Execute("₦ 1500", $false, $false, $false, $false, $false, `
    $true, 1, $false, "₦ 1500", 2)
